# #5: property aircraft done
#
# The "建物" (building) sheet rows were mistakenly tagged with
# property_category = "land"; fix them to "building".
# The "汽車" (car) sheet rows were mistakenly tagged with
# property_category = "land"; fix them to "car".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
for ($row = 2; $row -le 9; $row++) {
    $wsBuilding.Range("I" + $row).Value = "building"
}

$wsCar = $wb.Worksheets.Item("汽車")
for ($row = 2; $row -le 3; $row++) {
    $wsCar.Range("H" + $row).Value = "car"
}
